$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(765, 1).Value = 'TriggerAnnotation.Title'
$ws.Cells.Item(765, 2).Value = '触发注释'
$ws.Cells.Item(766, 1).Value = 'TriggerAnnotation.Text'
$ws.Cells.Item(766, 2).Value = '当前对象: %s'
$ws.Cells.Item(767, 1).Value = 'SaveMap_FileEncodingComment1_UTF8'
$ws.Cells.Item(767, 2).Value = '本文件编码为 UTF8，请使用此格式打开'
$ws.Cells.Item(768, 1).Value = 'Options.UTF8Support.InferEncoding'
$ws.Cells.Item(768, 2).Value = '打开ini和地图文件时自动推断编码'
$ws.Cells.Item(769, 1).Value = 'Options.UTF8Support.AlwaysSaveAsUTF8'
$ws.Cells.Item(769, 2).Value = '总是以UTF8编码保存地图'
$ws.Cells.Item(770, 1).Value = 'Menu.Display.Annotations'
$ws.Cells.Item(770, 2).Value = '地图注释\tAlt+9'
$ws.Cells.Item(771, 1).Value = 'Menu.Display.DamageFires'
$ws.Cells.Item(771, 2).Value = '受损火焰\tAlt+Num0'
$ws.Cells.Item(772, 1).Value = 'AllocFullMapBitmapFailed'
$ws.Cells.Item(772, 2).Value = '内存分配失败，无法渲染全图'
$ws.Cells.Item(773, 1).Value = 'MapRendererDlgCaption'
$ws.Cells.Item(773, 2).Value = '地图渲染选项'
$ws.Cells.Item(774, 1).Value = 'MapRendererDlgRenderSize'
$ws.Cells.Item(774, 2).Value = '渲染区域'
$ws.Cells.Item(775, 1).Value = 'MapRendererDlgRenderlayers'
$ws.Cells.Item(775, 2).Value = '渲染图层'
$ws.Cells.Item(776, 1).Value = 'MapRendererDlgLocalsize'
$ws.Cells.Item(776, 2).Value = '可见区域'
$ws.Cells.Item(777, 1).Value = 'MapRendererDlgFullsize'
$ws.Cells.Item(777, 2).Value = '全地图'
$ws.Cells.Item(778, 1).Value = 'MapRendererDlgIngame'
$ws.Cells.Item(778, 2).Value = '游戏内效果'
$ws.Cells.Item(779, 1).Value = 'MapRendererDlgCurrentlayers'
$ws.Cells.Item(779, 2).Value = '当前图层'
$ws.Cells.Item(780, 1).Value = 'MapRendererDlgTips'
$ws.Cells.Item(780, 2).Value = '地图渲染器会使用当前光照沙盒设置'
$ws.Cells.Item(781, 1).Value = 'Menu.Edit.TriggerAnnotation'
$ws.Cells.Item(781, 2).Value = '触发注释\tCtrl+Shift+A'
$ws.Cells.Item(782, 1).Value = 'Menu.MapTools.MapRenderer'
$ws.Cells.Item(782, 2).Value = '地图渲染器'
$ws.Cells.Item(783, 1).Value = 'MapRendererSuccess'
$ws.Cells.Item(783, 2).Value = '地图渲染已输出至：'

$ws.Range("B784").Select()

